$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 988.8
$ws.Range("I62").Value = 998
$ws.Range("J62").Value = 975
$ws.Range("K62").Value = 998
$ws.Range("L62").Value = 975
$ws.Range("M62").Value = -374
$ws.Range("N62").Value = -2223

$ws.Range("H64").Value = 3536.8076
$ws.Range("J64").Value = 3716.3076
$ws.Range("L64").Value = 3716.3076
$ws.Range("N64").Value = -4212.3076

$ws.Range("H65").Value = 988.8
$ws.Range("I65").Value = 998
$ws.Range("J65").Value = 975
$ws.Range("K65").Value = 4990
$ws.Range("L65").Value = 4875
$ws.Range("M65").Value = -1870
$ws.Range("N65").Value = -11115

$ws.Range("H67").Value = 3536.8076
$ws.Range("J67").Value = 3716.3076
$ws.Range("L67").Value = 3716.3076
$ws.Range("N67").Value = -5432.3076

$ws.Range("H111").Value = 1197.1666
$ws.Range("I111").Value = 988.25
$ws.Range("J111").Value = 1615
$ws.Range("K111").Value = 2964.75
$ws.Range("L111").Value = 4845
$ws.Range("M111").Value = 102.25
$ws.Range("N111").Value = -10979

$ws.Range("H121").Value = 851.6667
$ws.Range("J121").Value = 1000
$ws.Range("L121").Value = 3000
$ws.Range("N121").Value = -6494

$ws.Range("H125").Value = 2276.5
$ws.Range("I125").Value = 3061.5
$ws.Range("J125").Value = 1884
$ws.Range("K125").Value = 27553.5
$ws.Range("L125").Value = 16956
$ws.Range("M125").Value = -25093.5
$ws.Range("N125").Value = -21876

$ws.Range("H137").Value = 2734.6863
$ws.Range("I137").Value = 2501.658
$ws.Range("J137").Value = 3415.8462
$ws.Range("K137").Value = 7504.974
$ws.Range("L137").Value = 10247.5386
$ws.Range("M137").Value = -4954.974
$ws.Range("N137").Value = -15347.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6713.756
$ws.Range("I132").Value = 5841.8
$ws.Range("J132").Value = 8076.1875
$ws.Range("K132").Value = 17525.4
$ws.Range("L132").Value = 24228.5625
$ws.Range("M132").Value = -14995.4
$ws.Range("N132").Value = -29288.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value = 0

$ws.Range("H99").Value = 1658.7222
$ws.Range("I99").Value = 1576.25
$ws.Range("J99").Value = 1823.6666
$ws.Range("K99").Value = 1576.25
$ws.Range("L99").Value = 1823.6666
$ws.Range("M99").Value = -78.25
$ws.Range("N99").Value = -4819.6666

$ws.Range("H122").Value = 1963.875
$ws.Range("I122").Value = 1555.8334
$ws.Range("J122").Value = 3188
$ws.Range("K122").Value = 4667.5002
$ws.Range("L122").Value = 9564
$ws.Range("M122").Value = -2217.5002
$ws.Range("N122").Value = -14464

$ws.Range("H126").Value = 1658.7222
$ws.Range("I126").Value = 1576.25
$ws.Range("J126").Value = 1823.6666
$ws.Range("K126").Value = 4728.75
$ws.Range("L126").Value = 5470.9998
$ws.Range("M126").Value = -2258.75
$ws.Range("N126").Value = -10410.9998

$ws.Range("H132").Value = 17629.54
$ws.Range("I132").Value = 25753
$ws.Range("J132").Value = 14019.111
$ws.Range("K132").Value = 77259
$ws.Range("L132").Value = 42057.333
$ws.Range("M132").Value = -74729
$ws.Range("N132").Value = -47117.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 10000
$ws.Range("K42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("M42").Value = 30000
$ws.Range("N42").Value = -31068

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1597.7778
$ws.Range("I122").Value = 1395.6923
$ws.Range("J122").Value = 2123.2
$ws.Range("K122").Value = 4187.0769
$ws.Range("L122").Value = 6369.599999999999
$ws.Range("M122").Value = -1737.0769
$ws.Range("N122").Value = -11269.6

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2502
$ws.Range("I7").Value = 1200
$ws.Range("J7").Value = 2827.5
$ws.Range("K7").Value = 1200
$ws.Range("L7").Value = 2827.5
$ws.Range("M7").Value = -1088
$ws.Range("N7").Value = -3051.5

$ws.Range("H40").Value = 5556.5
$ws.Range("I40").Value = 5100.4443
$ws.Range("J40").Value = 6142.857
$ws.Range("K40").Value = 5100.4443
$ws.Range("L40").Value = 6142.857
$ws.Range("M40").Value = -4964.4443
$ws.Range("N40").Value = -6414.857

$ws.Range("H126").Value = 2502
$ws.Range("I126").Value = 1200
$ws.Range("J126").Value = 2827.5
$ws.Range("K126").Value = 3600
$ws.Range("L126").Value = 8482.5
$ws.Range("M126").Value = -1130
$ws.Range("N126").Value = -13422.5

$ws.Range("H132").Value = 5293.325
$ws.Range("I132").Value = 5922.3213
$ws.Range("J132").Value = 3825.6667
$ws.Range("K132").Value = 17766.9639
$ws.Range("L132").Value = 11477.0001
$ws.Range("M132").Value = -15236.9639
$ws.Range("N132").Value = -16537.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 32400
$ws.Range("I29").Value = 50000
$ws.Range("J29").Value = 28000
$ws.Range("K29").Value = 50000
$ws.Range("L29").Value = 28000
$ws.Range("M29").Value = -49710
$ws.Range("N29").Value = -28580

$ws.Range("H76").Value = 20150
$ws.Range("J76").Value = 20150
$ws.Range("L76").Value = 20150
$ws.Range("N76").Value = -20780

$ws.Range("H79").Value = 20150
$ws.Range("J79").Value = 20150
$ws.Range("L79").Value = 20150
$ws.Range("N79").Value = -22334

$ws.Range("H126").Value = 1305937.8
$ws.Range("I126").Value = 1114791.5
$ws.Range("J126").Value = 1673526.5
$ws.Range("K126").Value = 3344374.5
$ws.Range("L126").Value = 5020579.5
$ws.Range("M126").Value = -3341904.5
$ws.Range("N126").Value = -5025519.5

$ws.Range("H132").Value = 6229.6816
$ws.Range("I132").Value = 7291.4707
$ws.Range("J132").Value = 2619.6
$ws.Range("K132").Value = 21874.4121
$ws.Range("L132").Value = 7858.799999999999
$ws.Range("M132").Value = -19344.4121
$ws.Range("N132").Value = -12918.8
